$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("G2").Value = 66.47695399999999
$ws.Range("H2").Value = 199.430862
$ws.Range("I2").Value = 0.04311983106164722
$ws.Range("J2").Value = 0.04311983106164721
$ws.Range("O2").Value = 0.05546670559109387
$ws.Range("P2").Value = 0.05546670559109387
$ws.Range("Q2").Value = 4.075436141924
$ws.Range("R2").Value = 36.678925277316
$ws.Range("S2").Value = 0.002391714974634091
$ws.Range("T2").Value = 0.002391714974634091

$ws.Range("G3").Value = 66.47695399999999
$ws.Range("H3").Value = 199.430862
$ws.Range("I3").Value = 0.04311983106164722
$ws.Range("J3").Value = 0.04311983106164721
$ws.Range("M3").Value = 0.5397903333333334
$ws.Range("N3").Value = 1.619371
$ws.Range("O3").Value = 0.488376202980433
$ws.Range("P3").Value = 0.4883762029804329
$ws.Range("Q3").Value = 35.88361715864467
$ws.Range("R3").Value = 322.952554427802
$ws.Range("S3").Value = 0.02105869936704501
$ws.Range("T3").Value = 0.021058699367045

$ws.Range("G4").Value = 66.47695399999999
$ws.Range("H4").Value = 199.430862
$ws.Range("I4").Value = 0.04311983106164722
$ws.Range("J4").Value = 0.04311983106164721
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.005069
$ws.Range("N4").Value = 0.015207
$ws.Range("O4").Value = 0.004586186191257867
$ws.Range("P4").Value = 0.004586186191257867
$ws.Range("Q4").Value = 0.336971679826
$ws.Range("R4").Value = 3.032745118434
$ws.Range("S4").Value = 0.0001977555737842985
$ws.Range("T4").Value = 0.0001977555737842985

$ws.Range("G5").Value = 66.47695399999999
$ws.Range("H5").Value = 199.430862
$ws.Range("I5").Value = 0.04311983106164722
$ws.Range("J5").Value = 0.04311983106164721
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.4991103333333333
$ws.Range("N5").Value = 1.497331
$ws.Range("O5").Value = 0.4515709052372154
$ws.Range("P5").Value = 0.4515709052372153
$ws.Range("Q5").Value = 33.17933466992466
$ws.Range("R5").Value = 298.614012029322
$ws.Range("S5").Value = 0.01947166114618383
$ws.Range("T5").Value = 0.01947166114618383

$ws.Range("I6").Value = 0.8830494168872806
$ws.Range("J6").Value = 0.8830494168872804
$ws.Range("O6").Value = 0.05546670559109387
$ws.Range("P6").Value = 0.05546670559109387
$ws.Range("S6").Value = 0.04897984202887391
$ws.Range("T6").Value = 0.0489798420288739

$ws.Range("I7").Value = 0.8830494168872806
$ws.Range("J7").Value = 0.8830494168872804
$ws.Range("M7").Value = 0.5397903333333334
$ws.Range("N7").Value = 1.619371
$ws.Range("O7").Value = 0.488376202980433
$ws.Range("P7").Value = 0.4883762029804329
$ws.Range("Q7").Value = 734.8592614485329
$ws.Range("R7").Value = 6613.733353036797
$ws.Range("S7").Value = 0.4312603212634956
$ws.Range("T7").Value = 0.4312603212634954

$ws.Range("I8").Value = 0.8830494168872806
$ws.Range("J8").Value = 0.8830494168872804
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.005069
$ws.Range("N8").Value = 0.015207
$ws.Range("O8").Value = 0.004586186191257867
$ws.Range("P8").Value = 0.004586186191257867
$ws.Range("Q8").Value = 6.900830500760999
$ws.Range("R8").Value = 62.10747450684899
$ws.Range("S8").Value = 0.004049829041926758
$ws.Range("T8").Value = 0.004049829041926757

$ws.Range("I9").Value = 0.8830494168872806
$ws.Range("J9").Value = 0.8830494168872804
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.4991103333333333
$ws.Range("N9").Value = 1.497331
$ws.Range("O9").Value = 0.4515709052372154
$ws.Range("P9").Value = 0.4515709052372153
$ws.Range("Q9").Value = 679.4783609216129
$ws.Range("R9").Value = 6115.305248294517
$ws.Range("S9").Value = 0.3987594245529845
$ws.Range("T9").Value = 0.3987594245529844

$ws.Range("G10").Value = 44.831112
$ws.Range("H10").Value = 134.493336
$ws.Range("I10").Value = 0.02907940059566787
$ws.Range("J10").Value = 0.02907940059566786
$ws.Range("O10").Value = 0.05546670559109387
$ws.Range("P10").Value = 0.05546670559109387
$ws.Range("Q10").Value = 2.748416152272
$ws.Range("R10").Value = 24.735745370448
$ws.Range("S10").Value = 0.00161293855160539
$ws.Range("T10").Value = 0.001612938551605389

$ws.Range("G11").Value = 44.831112
$ws.Range("H11").Value = 134.493336
$ws.Range("I11").Value = 0.02907940059566787
$ws.Range("J11").Value = 0.02907940059566786
$ws.Range("M11").Value = 0.5397903333333334
$ws.Range("N11").Value = 1.619371
$ws.Range("O11").Value = 0.488376202980433
$ws.Range("P11").Value = 0.4883762029804329
$ws.Range("Q11").Value = 24.199400890184
$ws.Range("R11").Value = 217.794608011656
$ws.Range("S11").Value = 0.01420168724785922
$ws.Range("T11").Value = 0.01420168724785921

$ws.Range("G12").Value = 44.831112
$ws.Range("H12").Value = 134.493336
$ws.Range("I12").Value = 0.02907940059566787
$ws.Range("J12").Value = 0.02907940059566786
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.005069
$ws.Range("N12").Value = 0.015207
$ws.Range("O12").Value = 0.004586186191257867
$ws.Range("P12").Value = 0.004586186191257867
$ws.Range("Q12").Value = 0.227248906728
$ws.Range("R12").Value = 2.045240160552
$ws.Range("S12").Value = 0.0001333635454619078
$ws.Range("T12").Value = 0.0001333635454619077

$ws.Range("G13").Value = 44.831112
$ws.Range("H13").Value = 134.493336
$ws.Range("I13").Value = 0.02907940059566787
$ws.Range("J13").Value = 0.02907940059566786
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.4991103333333333
$ws.Range("N13").Value = 1.497331
$ws.Range("O13").Value = 0.4515709052372154
$ws.Range("P13").Value = 0.4515709052372153
$ws.Range("Q13").Value = 22.375671254024
$ws.Range("R13").Value = 201.381041286216
$ws.Range("S13").Value = 0.01313141125074136
$ws.Range("T13").Value = 0.01313141125074136

$ws.Range("G14").Value = 52.83062100000001
$ws.Range("H14").Value = 158.491863
$ws.Range("I14").Value = 0.0342682285413064
$ws.Range("J14").Value = 0.03426822854130639
$ws.Range("O14").Value = 0.05546670559109387
$ws.Range("P14").Value = 0.05546670559109387
$ws.Range("Q14").Value = 3.238834051026001
$ws.Range("R14").Value = 29.149506459234
$ws.Range("S14").Value = 0.001900745743628962
$ws.Range("T14").Value = 0.001900745743628962

$ws.Range("G15").Value = 52.83062100000001
$ws.Range("H15").Value = 158.491863
$ws.Range("I15").Value = 0.0342682285413064
$ws.Range("J15").Value = 0.03426822854130639
$ws.Range("M15").Value = 0.5397903333333334
$ws.Range("N15").Value = 1.619371
$ws.Range("O15").Value = 0.488376202980433
$ws.Range("P15").Value = 0.4883762029804329
$ws.Range("Q15").Value = 28.51745851979701
$ws.Range("R15").Value = 256.657126678173
$ws.Range("S15").Value = 0.01673578733786892
$ws.Range("T15").Value = 0.01673578733786892

$ws.Range("G16").Value = 52.83062100000001
$ws.Range("H16").Value = 158.491863
$ws.Range("I16").Value = 0.0342682285413064
$ws.Range("J16").Value = 0.03426822854130639
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.005069
$ws.Range("N16").Value = 0.015207
$ws.Range("O16").Value = 0.004586186191257867
$ws.Range("P16").Value = 0.004586186191257867
$ws.Range("Q16").Value = 0.267798417849
$ws.Range("R16").Value = 2.410185760641
$ws.Range("S16").Value = 0.0001571604765350081
$ws.Range("T16").Value = 0.0001571604765350081

$ws.Range("G17").Value = 52.83062100000001
$ws.Range("H17").Value = 158.491863
$ws.Range("I17").Value = 0.0342682285413064
$ws.Range("J17").Value = 0.03426822854130639
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.4991103333333333
$ws.Range("N17").Value = 1.497331
$ws.Range("O17").Value = 0.4515709052372154
$ws.Range("P17").Value = 0.4515709052372153
$ws.Range("Q17").Value = 26.368308857517
$ws.Range("R17").Value = 237.314779717653
$ws.Range("S17").Value = 0.01547453498327351
$ws.Range("T17").Value = 0.01547453498327351

$ws.Range("G18").Value = 16.16161433333333
$ws.Range("H18").Value = 48.484843
$ws.Range("I18").Value = 0.01048312291409786
$ws.Range("J18").Value = 0.01048312291409786
$ws.Range("O18").Value = 0.05546670559109387
$ws.Range("P18").Value = 0.05546670559109387
$ws.Range("Q18").Value = 0.9908039283193333
$ws.Range("R18").Value = 8.917235354874
$ws.Range("S18").Value = 0.0005814642923515163
$ws.Range("T18").Value = 0.0005814642923515161

$ws.Range("G19").Value = 16.16161433333333
$ws.Range("H19").Value = 48.484843
$ws.Range("I19").Value = 0.01048312291409786
$ws.Range("J19").Value = 0.01048312291409786
$ws.Range("M19").Value = 0.5397903333333334
$ws.Range("N19").Value = 1.619371
$ws.Range("O19").Value = 0.488376202980433
$ws.Range("P19").Value = 0.4883762029804329
$ws.Range("Q19").Value = 8.723883188194778
$ws.Range("R19").Value = 78.514948693753
$ws.Range("S19").Value = 0.005119707764164288
$ws.Range("T19").Value = 0.005119707764164285

$ws.Range("G20").Value = 16.16161433333333
$ws.Range("H20").Value = 48.484843
$ws.Range("I20").Value = 0.01048312291409786
$ws.Range("J20").Value = 0.01048312291409786
$ws.Range("K20").Value = 1
$ws.Range("L20").Value = 0.3333333333333333
$ws.Range("M20").Value = 0.005069
$ws.Range("N20").Value = 0.015207
$ws.Range("O20").Value = 0.004586186191257867
$ws.Range("P20").Value = 0.004586186191257867
$ws.Range("Q20").Value = 0.08192322305566666
$ws.Range("R20").Value = 0.7373090075009999
$ws.Range("S20").Value = 0.00004807755354989456
$ws.Range("T20").Value = 0.00004807755354989455

$ws.Range("G21").Value = 16.16161433333333
$ws.Range("H21").Value = 48.484843
$ws.Range("I21").Value = 0.01048312291409786
$ws.Range("J21").Value = 0.01048312291409786
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 0.4991103333333333
$ws.Range("N21").Value = 1.497331
$ws.Range("O21").Value = 0.4515709052372154
$ws.Range("P21").Value = 0.4515709052372153
$ws.Range("Q21").Value = 8.066428717114777
$ws.Range("R21").Value = 72.597858454033
$ws.Range("S21").Value = 0.004733873304032168
$ws.Range("T21").Value = 0.004733873304032166
